# Small SNPs.xlsx — rename sheets, repopulate "Random #2" (Sheet3) with a
# 49-sample / 5-SNP dataset, and update the view state (active tab, scroll
# position, selections) to match the post-edit workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Rename the sheets --------------------------------------------------
$ws1.Name = "All Samples, 5 SNPs"
$ws2.Name = "49 Samples, 5 SNPs"
$ws3.Name = "Random #2"

# --- Populate "Random #2" (Sheet3) with the 49-sample / 5-SNP table -----
$headers = @("rs3094315", "rs4475691", "rs3748597", "rs13303118", "rs9777703", "Phenotype")
for ($c = 0; $c -lt 6; $c++) {
    $ws3.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @(0,0,0,1,0,0),
    @(0,0,0,1,0,1),
    @(1,0,0,0,0,1),
    @(0,0,0,0,0,1),
    @(1,0,0,1,0,0),
    @(0,1,0,2,1,0),
    @(1,0,0,0,0,0),
    @(1,0,0,2,0,1),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,0,0,0),
    @(1,0,0,1,0,0),
    @(0,0,0,2,0,1),
    @(0,0,0,1,0,0),
    @(0,1,0,0,0,1),
    @(2,1,0,1,0,1),
    @(0,1,0,0,0,1),
    @(0,0,0,2,0,1),
    @(0,0,0,0,0,0),
    @(0,1,0,0,0,1),
    @(0,0,0,1,0,0),
    @(1,0,0,0,0,1),
    @(1,0,1,1,0,0),
    @(0,1,0,1,0,1),
    @(0,0,0,1,0,0),
    @(0,1,1,2,0,1),
    @(0,0,0,0,0,0),
    @(0,0,0,2,0,1),
    @(1,0,0,0,0,1),
    @(1,1,0,2,0,1),
    @(0,1,0,1,0,0),
    @(0,1,0,2,1,1),
    @(1,0,0,1,0,1),
    @(0,1,0,1,0,1),
    @(0,2,0,1,0,1),
    @(0,1,0,1,1,1),
    @(0,0,0,1,0,1),
    @(1,1,0,1,0,1),
    @(1,0,1,1,1,0),
    @(0,1,0,1,0,1),
    @(0,2,0,1,0,1),
    @(0,0,0,2,0,0),
    @(1,0,0,2,0,0),
    @(0,1,0,2,0,1),
    @(0,1,0,1,0,0),
    @(0,0,0,1,0,1),
    @(0,0,0,2,0,1),
    @(0,0,0,0,0,0),
    @(0,1,1,1,0,1)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws3.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- View state -----------------------------------------------------------
# Sheet1: scroll so row 34 is the top-left visible row, selection A1:F50
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A1:F50").Select()

# Sheet2: selection on F26 (no longer the tab-selected sheet)
$ws2.Activate()
$ws2.Range("F26").Select()

# Sheet3 ("Random #2") ends up the active / tab-selected sheet, selection H26
$ws3.Activate()
$ws3.Range("H26").Select()
